$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Insert a new worksheet named "CharChanges" right before the "Tiles" sheet
# ---------------------------------------------------------------------------
$tilesSheet = $wb.Worksheets.Item("Tiles")
$ws = $wb.Worksheets.Add($tilesSheet)
$ws.Name = "CharChanges"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.16
$ws.Columns.Item(2).ColumnWidth = 45.42

# ---------------------------------------------------------------------------
# Header row (bold, yellow fill, boxed border)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Char"
$ws.Range("B1").Value = "Change"

$headerRng = $ws.Range("A1:B1")
$headerRng.Font.Bold = $true
$headerRng.Interior.Color = 65535
$headerRng.Rows.Item(1).RowHeight = 15.75

$headerRng.Borders.Item(8).Weight = -4138
$headerRng.Borders.Item(9).Weight = -4138
$headerRng.Borders.Item(7).Weight = -4138
$headerRng.Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Leonaria"
$ws.Range("B2").Value = "SLP/Lvl 10 -> 20"

$ws.Range("A3").Value = "Leonaria"
$ws.Range("B3").Value = "Changed start spells"

$ws.Range("A4").Value = "Targor"
$ws.Range("B4").Value = "Start SLP 20 -> 25"

$ws.Range("A5").Value = "Targor"
$ws.Range("B5").Value = "Changed start spells"

$ws.Range("A6").Value = "Leonaria"
$ws.Range("B6").Value = "Replace some start items"

$ws.Range("A7").Value = "Targor"
$ws.Range("B7").Value = "Replace some start items"

$ws.Range("A8").Value = "Valdyn"
$ws.Range("B8").Value = "Added Monster Knowledge spell scroll"

$ws.Range("A9").Value = "Valdyn"
$ws.Range("B9").Value = "Start SLP 16 -> 10"

$ws.Range("A10").Value = "Valdyn"
$ws.Range("B10").Value = "Remove Monster Knowledge as start spell"

# ---------------------------------------------------------------------------
# Make this the active sheet / selection, matching the saved view state
# ---------------------------------------------------------------------------
$ws.Range("E13").Select()
